$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "取得日時" (retrieved at) timestamp for the newly appended
# batch of rows (2-11) to reflect the latest scrape run.
$newTimestamp = "2025-11-17 18:33:23"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
